# Update the PSSM score grid (B2:K21) with recomputed values
# per the supplemental-figures data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = -20.70199593005316
$ws.Cells.Item(2, 3).Value = 2.344629360788994
$ws.Cells.Item(2, 4).Value = -20.70199593005316
$ws.Cells.Item(2, 5).Value = -20.70199593005316
$ws.Cells.Item(2, 6).Value = -20.70199593005316
$ws.Cells.Item(2, 7).Value = -20.70199593005316
$ws.Cells.Item(2, 8).Value = -20.70199593005316
$ws.Cells.Item(2, 9).Value = -20.70199593005316
$ws.Cells.Item(2, 10).Value = -20.70199593005316
$ws.Cells.Item(2, 11).Value = -20.70199593005316

# Row 3
$ws.Cells.Item(3, 2).Value = -20.70199593005316
$ws.Cells.Item(3, 3).Value = -20.70199593005316
$ws.Cells.Item(3, 4).Value = -20.70199593005316
$ws.Cells.Item(3, 5).Value = -20.70199593005316
$ws.Cells.Item(3, 6).Value = -20.70199593005316
$ws.Cells.Item(3, 7).Value = -20.70199593005316
$ws.Cells.Item(3, 8).Value = -20.70199593005316
$ws.Cells.Item(3, 9).Value = 2.389754031364213
$ws.Cells.Item(3, 10).Value = -20.70199593005316
$ws.Cells.Item(3, 11).Value = -20.70199593005316

# Row 4
$ws.Cells.Item(4, 2).Value = -20.70199593005316
$ws.Cells.Item(4, 3).Value = 2.258997430661782
$ws.Cells.Item(4, 4).Value = 2.778320088105378
$ws.Cells.Item(4, 5).Value = -20.70199593005316
$ws.Cells.Item(4, 6).Value = 2.568906884459628
$ws.Cells.Item(4, 7).Value = -20.70199593005316
$ws.Cells.Item(4, 8).Value = 1.832827707688418
$ws.Cells.Item(4, 9).Value = -20.70199593005316
$ws.Cells.Item(4, 10).Value = 2.389773397522596
$ws.Cells.Item(4, 11).Value = -20.70199593005316

# Row 5
$ws.Cells.Item(5, 2).Value = -20.70199593005316
$ws.Cells.Item(5, 3).Value = 1.056904316413608
$ws.Cells.Item(5, 4).Value = -20.70199593005316
$ws.Cells.Item(5, 5).Value = -20.70199593005316
$ws.Cells.Item(5, 6).Value = -20.70199593005316
$ws.Cells.Item(5, 7).Value = 2.098066018351406
$ws.Cells.Item(5, 8).Value = -20.70199593005316
$ws.Cells.Item(5, 9).Value = -20.70199593005316
$ws.Cells.Item(5, 10).Value = -20.70199593005316
$ws.Cells.Item(5, 11).Value = -20.70199593005316

# Row 6
$ws.Cells.Item(6, 2).Value = -20.70199593005316
$ws.Cells.Item(6, 3).Value = -20.70199593005316
$ws.Cells.Item(6, 4).Value = -20.70199593005316
$ws.Cells.Item(6, 5).Value = -20.70199593005316
$ws.Cells.Item(6, 6).Value = -20.70199593005316
$ws.Cells.Item(6, 7).Value = -20.70199593005316
$ws.Cells.Item(6, 8).Value = -20.70199593005316
$ws.Cells.Item(6, 9).Value = -20.70199593005316
$ws.Cells.Item(6, 10).Value = -20.70199593005316
$ws.Cells.Item(6, 11).Value = -20.70199593005316

# Row 7
$ws.Cells.Item(7, 2).Value = -20.70199593005316
$ws.Cells.Item(7, 3).Value = -20.70199593005316
$ws.Cells.Item(7, 4).Value = -20.70199593005316
$ws.Cells.Item(7, 5).Value = -20.70199593005316
$ws.Cells.Item(7, 6).Value = -20.70199593005316
$ws.Cells.Item(7, 7).Value = -20.70199593005316
$ws.Cells.Item(7, 8).Value = -20.70199593005316
$ws.Cells.Item(7, 9).Value = -20.70199593005316
$ws.Cells.Item(7, 10).Value = -20.70199593005316
$ws.Cells.Item(7, 11).Value = -20.70199593005316

# Row 8
$ws.Cells.Item(8, 2).Value = -20.70199593005316
$ws.Cells.Item(8, 3).Value = -20.70199593005316
$ws.Cells.Item(8, 4).Value = -20.70199593005316
$ws.Cells.Item(8, 5).Value = 2.849330845843308
$ws.Cells.Item(8, 6).Value = -20.70199593005316
$ws.Cells.Item(8, 7).Value = -20.70199593005316
$ws.Cells.Item(8, 8).Value = -20.70199593005316
$ws.Cells.Item(8, 9).Value = -20.70199593005316
$ws.Cells.Item(8, 10).Value = -20.70199593005316
$ws.Cells.Item(8, 11).Value = -20.70199593005316

# Row 9
$ws.Cells.Item(9, 2).Value = 4.321927291404681
$ws.Cells.Item(9, 3).Value = -20.70199593005316
$ws.Cells.Item(9, 4).Value = -20.70199593005316
$ws.Cells.Item(9, 5).Value = -20.70199593005316
$ws.Cells.Item(9, 6).Value = -20.70199593005316
$ws.Cells.Item(9, 7).Value = -20.70199593005316
$ws.Cells.Item(9, 8).Value = -20.70199593005316
$ws.Cells.Item(9, 9).Value = -20.70199593005316
$ws.Cells.Item(9, 10).Value = -20.70199593005316
$ws.Cells.Item(9, 11).Value = -20.70199593005316

# Row 10
$ws.Cells.Item(10, 2).Value = -20.70199593005316
$ws.Cells.Item(10, 3).Value = -20.70199593005316
$ws.Cells.Item(10, 4).Value = -20.70199593005316
$ws.Cells.Item(10, 5).Value = -20.70199593005316
$ws.Cells.Item(10, 6).Value = -20.70199593005316
$ws.Cells.Item(10, 7).Value = -20.70199593005316
$ws.Cells.Item(10, 8).Value = -20.70199593005316
$ws.Cells.Item(10, 9).Value = 1.537817294953475
$ws.Cells.Item(10, 10).Value = -20.70199593005316
$ws.Cells.Item(10, 11).Value = 2.244228045877663

# Row 11
$ws.Cells.Item(11, 2).Value = -20.70199593005316
$ws.Cells.Item(11, 3).Value = -20.70199593005316
$ws.Cells.Item(11, 4).Value = -20.70199593005316
$ws.Cells.Item(11, 5).Value = 2.004003605698548
$ws.Cells.Item(11, 6).Value = -20.70199593005316
$ws.Cells.Item(11, 7).Value = 2.577337636869133
$ws.Cells.Item(11, 8).Value = -20.70199593005316
$ws.Cells.Item(11, 9).Value = -20.70199593005316
$ws.Cells.Item(11, 10).Value = -20.70199593005316
$ws.Cells.Item(11, 11).Value = 1.381185034068269

# Row 12
$ws.Cells.Item(12, 2).Value = -20.70199593005316
$ws.Cells.Item(12, 3).Value = -20.70199593005316
$ws.Cells.Item(12, 4).Value = -20.70199593005316
$ws.Cells.Item(12, 5).Value = -20.70199593005316
$ws.Cells.Item(12, 6).Value = -20.70199593005316
$ws.Cells.Item(12, 7).Value = -20.70199593005316
$ws.Cells.Item(12, 8).Value = -20.70199593005316
$ws.Cells.Item(12, 9).Value = -20.70199593005316
$ws.Cells.Item(12, 10).Value = -20.70199593005316
$ws.Cells.Item(12, 11).Value = -20.70199593005316

# Row 13
$ws.Cells.Item(13, 2).Value = -20.70199593005316
$ws.Cells.Item(13, 3).Value = -20.70199593005316
$ws.Cells.Item(13, 4).Value = -20.70199593005316
$ws.Cells.Item(13, 5).Value = 1.726684477243466
$ws.Cells.Item(13, 6).Value = -20.70199593005316
$ws.Cells.Item(13, 7).Value = -20.70199593005316
$ws.Cells.Item(13, 8).Value = -20.70199593005316
$ws.Cells.Item(13, 9).Value = -20.70199593005316
$ws.Cells.Item(13, 10).Value = 2.276947171288418
$ws.Cells.Item(13, 11).Value = 1.608323119756478

# Row 14
$ws.Cells.Item(14, 2).Value = -20.70199593005316
$ws.Cells.Item(14, 3).Value = -20.70199593005316
$ws.Cells.Item(14, 4).Value = 1.651563645248576
$ws.Cells.Item(14, 5).Value = -20.70199593005316
$ws.Cells.Item(14, 6).Value = -20.70199593005316
$ws.Cells.Item(14, 7).Value = -20.70199593005316
$ws.Cells.Item(14, 8).Value = -20.70199593005316
$ws.Cells.Item(14, 9).Value = -20.70199593005316
$ws.Cells.Item(14, 10).Value = -20.70199593005316
$ws.Cells.Item(14, 11).Value = 2.076611639052086

# Row 15
$ws.Cells.Item(15, 2).Value = -20.70199593005316
$ws.Cells.Item(15, 3).Value = -20.70199593005316
$ws.Cells.Item(15, 4).Value = -0.2873509485754108
$ws.Cells.Item(15, 5).Value = -20.70199593005316
$ws.Cells.Item(15, 6).Value = -20.70199593005316
$ws.Cells.Item(15, 7).Value = -20.70199593005316
$ws.Cells.Item(15, 8).Value = -20.70199593005316
$ws.Cells.Item(15, 9).Value = -20.70199593005316
$ws.Cells.Item(15, 10).Value = -20.70199593005316
$ws.Cells.Item(15, 11).Value = -20.70199593005316

# Row 16
$ws.Cells.Item(16, 2).Value = -20.70199593005316
$ws.Cells.Item(16, 3).Value = -20.70199593005316
$ws.Cells.Item(16, 4).Value = -20.70199593005316
$ws.Cells.Item(16, 5).Value = -20.70199593005316
$ws.Cells.Item(16, 6).Value = -20.70199593005316
$ws.Cells.Item(16, 7).Value = -20.70199593005316
$ws.Cells.Item(16, 8).Value = -20.70199593005316
$ws.Cells.Item(16, 9).Value = -20.70199593005316
$ws.Cells.Item(16, 10).Value = 2.30454399910942
$ws.Cells.Item(16, 11).Value = -20.70199593005316

# Row 17
$ws.Cells.Item(17, 2).Value = -20.70199593005316
$ws.Cells.Item(17, 3).Value = 0.7304771557224392
$ws.Cells.Item(17, 4).Value = 0.1452855041657443
$ws.Cells.Item(17, 5).Value = -20.70199593005316
$ws.Cells.Item(17, 6).Value = -20.70199593005316
$ws.Cells.Item(17, 7).Value = -20.70199593005316
$ws.Cells.Item(17, 8).Value = 0.5399133416149213
$ws.Cells.Item(17, 9).Value = 0.9125238642253763
$ws.Cells.Item(17, 10).Value = 1.252714771946028
$ws.Cells.Item(17, 11).Value = -20.70199593005316

# Row 18
$ws.Cells.Item(18, 2).Value = -20.70199593005316
$ws.Cells.Item(18, 3).Value = -20.70199593005316
$ws.Cells.Item(18, 4).Value = -20.70199593005316
$ws.Cells.Item(18, 5).Value = -20.70199593005316
$ws.Cells.Item(18, 6).Value = -20.70199593005316
$ws.Cells.Item(18, 7).Value = -20.70199593005316
$ws.Cells.Item(18, 8).Value = 0.4655453345497921
$ws.Cells.Item(18, 9).Value = 0.9417374176033035
$ws.Cells.Item(18, 10).Value = 1.372796568894604
$ws.Cells.Item(18, 11).Value = -20.70199593005316

# Row 19
$ws.Cells.Item(19, 2).Value = -20.70199593005316
$ws.Cells.Item(19, 3).Value = -20.70199593005316
$ws.Cells.Item(19, 4).Value = 1.811714560699863
$ws.Cells.Item(19, 5).Value = -20.70199593005316
$ws.Cells.Item(19, 6).Value = -20.70199593005316
$ws.Cells.Item(19, 7).Value = -20.70199593005316
$ws.Cells.Item(19, 8).Value = 1.884736905336426
$ws.Cells.Item(19, 9).Value = 2.087518014315785
$ws.Cells.Item(19, 10).Value = -20.70199593005316
$ws.Cells.Item(19, 11).Value = -20.70199593005316

# Row 20
$ws.Cells.Item(20, 2).Value = -20.70199593005316
$ws.Cells.Item(20, 3).Value = 1.62879870487825
$ws.Cells.Item(20, 4).Value = 2.189595151323978
$ws.Cells.Item(20, 5).Value = -20.70199593005316
$ws.Cells.Item(20, 6).Value = 3.814180493987279
$ws.Cells.Item(20, 7).Value = -20.70199593005316
$ws.Cells.Item(20, 8).Value = 2.196353954812346
$ws.Cells.Item(20, 9).Value = 1.926919679957833
$ws.Cells.Item(20, 10).Value = -20.70199593005316
$ws.Cells.Item(20, 11).Value = 2.430340780535921

# Row 21
$ws.Cells.Item(21, 2).Value = -20.70199593005316
$ws.Cells.Item(21, 3).Value = 1.723297075448558
$ws.Cells.Item(21, 4).Value = -20.70199593005316
$ws.Cells.Item(21, 5).Value = 2.452213543428289
$ws.Cells.Item(21, 6).Value = -20.70199593005316
$ws.Cells.Item(21, 7).Value = 3.285441665950212
$ws.Cells.Item(21, 8).Value = 2.413287965915229
$ws.Cells.Item(21, 9).Value = -20.70199593005316
$ws.Cells.Item(21, 10).Value = -20.70199593005316
$ws.Cells.Item(21, 11).Value = -20.70199593005316

